$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.212.49'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '1.855.84'
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''235.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '''0.4776'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.64%  '
$ws.Range("D8").Value = '''0.2800'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.21%  '
$ws.Range("D9").Value = '''0.06461'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.46%  '
$ws.Range("D10").Value = '1.854.87'
$ws.Range("E10").Value = '  -1.83%  '
$ws.Range("D11").Value = '''0.07358'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("E12").Value = '  -4.12%  '
$ws.Range("D13").Value = '''5.078'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("D14").Value = '''87.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = '''0.6448'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.00%  '
$ws.Range("D16").Value = '30.155.45'
$ws.Range("E16").Value = '  -0.87%  '
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '''13.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("D19").Value = '''0.000007573'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("D20").Value = '''225.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +16.90%  '
$ws.Range("D21").Value = '2.097.78'
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("D24").Value = '''6.074'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").Value = '''9.195'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.03%  '
$ws.Range("D26").Value = '''163.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").Value = '''18.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").Value = '''1.923'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '''0.09181'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''4.229'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("D33").Value = '''0.04959'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.04%  '
$ws.Range("D34").Value = '''0.7315'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("D35").Value = '''1.143'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.10%  '
$ws.Range("D37").Value = '''0.01838'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.98%  '
$ws.Range("D38").Value = '''2.592'
$ws.Range("D38").Style = "Normal"
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''2.048'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.58%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''0.8985'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("D41").Value = '''5.947'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("D42").Value = '''105.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").Value = '''1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = '''0.4226'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.45%  '
$ws.Range("D45").Value = '''7.349'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("E46").Value = '  -3.72%  '
$ws.Range("D47").Value = '''64.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.12%  '
$ws.Range("D48").Value = '''1.501'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.14%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.674'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.94%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''33.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.97%  '
$ws.Range("D51").Value = '''0.05651'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.24%  '
